$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I (I0) and J (IF), matching the style of the other header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

# Data rows 2-28: column I = 1 (constant), column J = same value as column H
for ($r = 2; $r -le 28; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
